$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix casing of the "RiboPure0.25x" -> "RiboPure0.25X" label used in G2:G13.
$ws.Range("G2:G13").Value = "RiboPure0.25X"

# 2. Column H (rows 2-13) used to hold a plain boolean FALSE constant; turn each
#    one into an explicit =FALSE() formula (set cell-by-cell so Excel doesn't
#    collapse them into a single shared formula spanning the range).
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# 3. Move the active selection from H2:H13 to the single cell G13.
$ws.Range("G13").Select() | Out-Null
